$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for Alpha / Adjusted R2, copying the header style (bold, border, centered)
$ws.Range("A1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "Alpha"
$ws.Range("J1").Value = "Adjusted R2"

# Row 2 (Accruals) - Newey-West model results
$ws.Range("C2").Value = 0.4847547169117808
$ws.Range("E2").Value = 0.6986757104210642
$ws.Range("I2").Value = -0.004598563622157524
$ws.Range("J2").Value = 0.1189339105006265

# Row 3 (Assest Growth)
$ws.Range("C3").Value = 0.4964198355513971
$ws.Range("E3").Value = 0.6801336561869071
$ws.Range("I3").Value = -0.003949309491163033
$ws.Range("J3").Value = 0.1518558948927919

# Row 4 (BM)
$ws.Range("C4").Value = 0.09464160838125377
$ws.Range("D4").Value = $false
$ws.Range("E4").Value = -1.671405509699841
$ws.Range("I4").Value = 0.01570082964427109
$ws.Range("J4").Value = 0.3354477810839157

# Row 5 (Gross Profit)
$ws.Range("C5").Value = 0.6170892749780951
$ws.Range("E5").Value = -0.499979836903159
$ws.Range("I5").Value = 0.006475692414618974
$ws.Range("J5").Value = 0.1969900145097194

# Row 6 (Momentum)
$ws.Range("C6").Value = 0.1994279935708438
$ws.Range("D6").Value = $false
$ws.Range("E6").Value = 1.28318293385164
$ws.Range("I6").Value = -0.03991517995090906
$ws.Range("J6").Value = 0.2188846263134701

# Row 7 (Leaverage Ret)
$ws.Range("C7").Value = 0.3432046942032331
$ws.Range("E7").Value = -0.9478522153205311
$ws.Range("I7").Value = 0.006393198039684209
$ws.Range("J7").Value = 0.1524569349572032
